# Apply "hybrid bold + color" highlighting to quantitative metrics
# (percentages, dollar amounts, large numbers) in specific bullet
# paragraphs, matching the target diff.
#
# Strategy: for each target paragraph, scan left-to-right with a scoped
# Find (Range restricted to that paragraph) for each metric substring in
# order, and apply Font.Bold + Font.Color (dark slate, hex 2C3E50) to the
# matched sub-range. Word's COM engine splits runs automatically when a
# sub-range's formatting is changed, which reproduces the <w:r> splitting
# seen in the diff.

$d = $word.ActiveDocument

# 2C3E50 (R=0x2C,G=0x3E,B=0x50) expressed as a Word "OLE" color long,
# which is encoded 0x00BBGGRR.
$metricColor = 5258796

function Set-MetricHighlight {
    param(
        [int]$ParagraphIndex,
        [string[]]$Metrics
    )

    $p = $d.Paragraphs($ParagraphIndex)
    $paraEnd = $p.Range.End
    $cursor = $p.Range.Start

    foreach ($metric in $Metrics) {
        $rng = $d.Range($cursor, $paraEnd)
        $rng.Find.Text = $metric
        $rng.Find.MatchCase = $true
        $found = $rng.Find.Execute()
        if ($found) {
            $rng.Font.Bold = 1
            $rng.Font.Color = $metricColor
            $cursor = $rng.End
        }
    }
}

# Paragraph 9 (1-based): "...demographic classification accuracy from 23% to 64%"
Set-MetricHighlight 9 @("23%", "64%")

# Paragraph 11 (1-based): "Achieved 87% ... 71% ... from ±4.2% to ±2.1%"
Set-MetricHighlight 11 @("87%", "71%", "±4.2%", "±2.1%")

# Paragraph 31 (1-based): "Wrote RFP and analyzed bids from 1,200 vendors..."
Set-MetricHighlight 31 @("1,200")

# Paragraph 46 (1-based): "...became the $400M Polling Consortium Database... now valued at $1B+"
Set-MetricHighlight 46 @('$400M', '$1B')

# Paragraph 63 (1-based): "Algorithm reduced mapping costs by 73.5%, ... $4.7M"
Set-MetricHighlight 63 @('73.5%', '$4.7M')

# Paragraph 65 (1-based): "Achieved 87% prediction accuracy ... industry standard of 71%"
Set-MetricHighlight 65 @("87%", "71%")

Write-Output "Metric highlighting applied."
